$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1000
$ws.Range("M32").Value = ""
$ws.Range("N32").Value = -1652
$ws.Range("H40").Value = 2447.4
$ws.Range("I40").Value = 3500.2
$ws.Range("J40").Value = 1921
$ws.Range("K40").Value = 3500.2
$ws.Range("L40").Value = 1921
$ws.Range("M40").Value = -3325.2
$ws.Range("N40").Value = -2271
$ws.Range("H45").Value = 27803528
$ws.Range("I45").Value = 5000
$ws.Range("J45").Value = 37069704
$ws.Range("K45").Value = 15000
$ws.Range("L45").Value = 111209112
$ws.Range("M45").Value = -14808
$ws.Range("N45").Value = -111209496
$ws.Range("H46").Value = 5116
$ws.Range("J46").Value = 4463.4287
$ws.Range("L46").Value = 13390.2861
$ws.Range("N46").Value = -13628.2861
$ws.Range("H49").Value = 19231072
$ws.Range("J49").Value = 76923070
$ws.Range("L49").Value = 230769210
$ws.Range("N49").Value = -230769482
$ws.Range("H59").Value = 2977306.5
$ws.Range("I59").Value = 908.5
$ws.Range("J59").Value = 4167865.5
$ws.Range("K59").Value = 2725.5
$ws.Range("L59").Value = 12503596.5
$ws.Range("M59").Value = -2168.5
$ws.Range("N59").Value = -12504710.5
$ws.Range("H60").Value = 5116
$ws.Range("J60").Value = 4463.4287
$ws.Range("L60").Value = 13390.2861
$ws.Range("N60").Value = -14358.2861
$ws.Range("H116").Value = 2081.818
$ws.Range("I116").Value = 2016.6666
$ws.Range("J116").Value = 2160
$ws.Range("K116").Value = 2016.6666
$ws.Range("L116").Value = 2160
$ws.Range("M116").Value = 1425.3334
$ws.Range("N116").Value = -9044
$ws.Range("H140").Value = 89962.5
$ws.Range("J140").Value = 89962.5
$ws.Range("L140").Value = 89962.5
$ws.Range("N140").Value = -100322.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1624.6111
$ws.Range("I2").Value = 1624.6111
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1624.6111
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1511.6111
$ws.Range("N2").Value = ""
$ws.Range("H45").Value = 1384
$ws.Range("I45").Value = 1415.1923
$ws.Range("J45").Value = 1248.8334
$ws.Range("K45").Value = 1415.1923
$ws.Range("L45").Value = 1248.8334
$ws.Range("M45").Value = -1038.1923
$ws.Range("N45").Value = -2002.8334
$ws.Range("H74").Value = 4180.6924
$ws.Range("I74").Value = 1655.4
$ws.Range("J74").Value = 26277
$ws.Range("K74").Value = 1655.4
$ws.Range("L74").Value = 26277
$ws.Range("M74").Value = -781.4000000000001
$ws.Range("N74").Value = -28025
$ws.Range("H77").Value = 4180.6924
$ws.Range("I77").Value = 1655.4
$ws.Range("J77").Value = 26277
$ws.Range("K77").Value = 8277
$ws.Range("L77").Value = 131385
$ws.Range("M77").Value = -3909
$ws.Range("N77").Value = -140121
$ws.Range("H116").Value = 1624.6111
$ws.Range("I116").Value = 1624.6111
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1624.6111
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 669.3888999999999
$ws.Range("N116").Value = ""

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1624.6111
$ws.Range("I3").Value = 1624.6111
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1624.6111
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1510.6111
$ws.Range("N3").Value = ""
$ws.Range("H80").Value = 234.8
$ws.Range("I80").Value = 147.5
$ws.Range("J80").Value = 248.23077
$ws.Range("K80").Value = 147.5
$ws.Range("L80").Value = 248.23077
$ws.Range("M80").Value = 850.5
$ws.Range("N80").Value = -2244.23077
$ws.Range("H83").Value = 234.8
$ws.Range("I83").Value = 147.5
$ws.Range("J83").Value = 248.23077
$ws.Range("K83").Value = 737.5
$ws.Range("L83").Value = 1241.15385
$ws.Range("M83").Value = 4254.5
$ws.Range("N83").Value = -11225.15385
$ws.Range("H99").Value = 1068.6666
$ws.Range("I99").Value = 874
$ws.Range("J99").Value = 1750
$ws.Range("K99").Value = 874
$ws.Range("L99").Value = 1750
$ws.Range("M99").Value = 624
$ws.Range("N99").Value = -4746

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3095
$ws.Range("I62").Value = 2875.8333
$ws.Range("J62").Value = 3533.3333
$ws.Range("K62").Value = 2875.8333
$ws.Range("L62").Value = 3533.3333
$ws.Range("M62").Value = -2251.8333
$ws.Range("N62").Value = -4781.3333
$ws.Range("H65").Value = 3095
$ws.Range("I65").Value = 2875.8333
$ws.Range("J65").Value = 3533.3333
$ws.Range("K65").Value = 14379.1665
$ws.Range("L65").Value = 17666.6665
$ws.Range("M65").Value = -11259.1665
$ws.Range("N65").Value = -23906.6665
$ws.Range("H99").Value = 1480.4615
$ws.Range("I99").Value = 1077.3334
$ws.Range("J99").Value = 2387.5
$ws.Range("K99").Value = 1077.3334
$ws.Range("L99").Value = 2387.5
$ws.Range("M99").Value = 420.6666
$ws.Range("N99").Value = -5383.5
$ws.Range("H105").Value = 623.2222
$ws.Range("I105").Value = 576.25
$ws.Range("J105").Value = 999
$ws.Range("K105").Value = 576.25
$ws.Range("L105").Value = 999
$ws.Range("M105").Value = 1170.75
$ws.Range("N105").Value = -4493
$ws.Range("H107").Value = 723.05554
$ws.Range("I107").Value = 411.83334
$ws.Range("J107").Value = 878.6667
$ws.Range("K107").Value = 411.83334
$ws.Range("L107").Value = 878.6667
$ws.Range("M107").Value = 1508.16666
$ws.Range("N107").Value = -4718.6667
$ws.Range("H122").Value = 5030
$ws.Range("I122").Value = 5304
$ws.Range("J122").Value = 920
$ws.Range("K122").Value = 15912
$ws.Range("L122").Value = 2760
$ws.Range("M122").Value = -13462
$ws.Range("N122").Value = -7660
$ws.Range("H126").Value = 1480.4615
$ws.Range("I126").Value = 1077.3334
$ws.Range("J126").Value = 2387.5
$ws.Range("K126").Value = 3232.0002
$ws.Range("L126").Value = 7162.5
$ws.Range("M126").Value = -762.0001999999999
$ws.Range("N126").Value = -12102.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 12858425
$ws.Range("J9").Value = 12858425
$ws.Range("L9").Value = 38575275
$ws.Range("N9").Value = -38575723
$ws.Range("H38").Value = 47.916668
$ws.Range("I38").Value = 30.555555
$ws.Range("J38").Value = 100
$ws.Range("K38").Value = 91.66666499999999
$ws.Range("L38").Value = 300
$ws.Range("M38").Value = 255.333335
$ws.Range("N38").Value = -994
$ws.Range("H76").Value = 3202.3333
$ws.Range("I76").Value = 2933
$ws.Range("J76").Value = 3471.6667
$ws.Range("K76").Value = 8799
$ws.Range("L76").Value = 10415.0001
$ws.Range("M76").Value = -8416
$ws.Range("N76").Value = -11181.0001
$ws.Range("H79").Value = 3202.3333
$ws.Range("I79").Value = 2933
$ws.Range("J79").Value = 3471.6667
$ws.Range("K79").Value = 8799
$ws.Range("L79").Value = 10415.0001
$ws.Range("M79").Value = -7473
$ws.Range("N79").Value = -13067.0001
$ws.Range("H108").Value = 976.8570999999999
$ws.Range("I108").Value = 976.8570999999999
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 2930.5713
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -50.57129999999961
$ws.Range("N108").Value = ""

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5957.4165
$ws.Range("I122").Value = 6944.3335
$ws.Range("J122").Value = 2996.6667
$ws.Range("K122").Value = 20833.0005
$ws.Range("L122").Value = 8990.000100000001
$ws.Range("M122").Value = -18383.0005
$ws.Range("N122").Value = -13890.0001
$ws.Range("H126").Value = 2354.913
$ws.Range("I126").Value = 1816.5834
$ws.Range("J126").Value = 2942.182
$ws.Range("K126").Value = 5449.7502
$ws.Range("L126").Value = 8826.545999999998
$ws.Range("M126").Value = -2979.7502
$ws.Range("N126").Value = -13766.546
$ws.Range("H141").Value = 51151.6
$ws.Range("J141").Value = 51151.6
$ws.Range("L141").Value = 51151.6
$ws.Range("N141").Value = -61511.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3770.647
$ws.Range("I7").Value = 4230.3
$ws.Range("J7").Value = 3114
$ws.Range("K7").Value = 4230.3
$ws.Range("L7").Value = 3114
$ws.Range("M7").Value = -4118.3
$ws.Range("N7").Value = -3338
$ws.Range("H93").Value = 991.4
$ws.Range("I93").Value = 788.5
$ws.Range("J93").Value = 1126.6666
$ws.Range("K93").Value = 788.5
$ws.Range("L93").Value = 1126.6666
$ws.Range("M93").Value = 459.5
$ws.Range("N93").Value = -3622.6666
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = ""
$ws.Range("H126").Value = 3770.647
$ws.Range("I126").Value = 4230.3
$ws.Range("J126").Value = 3114
$ws.Range("K126").Value = 12690.9
$ws.Range("L126").Value = 9342
$ws.Range("M126").Value = -10220.9
$ws.Range("N126").Value = -14282

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 49802
$ws.Range("J103").Value = 49802
$ws.Range("L103").Value = 49802
$ws.Range("N103").Value = -52146
$ws.Range("H107").Value = 2499.2778
$ws.Range("I107").Value = 344
$ws.Range("J107").Value = 3328.2307
$ws.Range("K107").Value = 1032
$ws.Range("L107").Value = 9984.6921
$ws.Range("M107").Value = 888
$ws.Range("N107").Value = -13824.6921
$ws.Range("H108").Value = 21000
$ws.Range("J108").Value = 21000
$ws.Range("L108").Value = 21000
$ws.Range("N108").Value = -28680
$ws.Range("H140").Value = 59761.285
$ws.Range("J140").Value = 59761.285
$ws.Range("L140").Value = 59761.285
$ws.Range("N140").Value = -70121.285
$ws.Range("H141").Value = 40580
$ws.Range("J141").Value = 40580
$ws.Range("L141").Value = 40580
$ws.Range("N141").Value = -50940
